$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3..56 down to 4..57)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly record
$ws.Range("A3").Value = 9
$ws.Range("B3").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C3").Value = 'Metropolitana'
$ws.Range("D3").Value = '2021-10-21'
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 'Fruta'
$ws.Range("G3").Value = 100102
$ws.Range("H3").Value = 'Cítricos'
$ws.Range("I3").Value = 100102006
$ws.Range("J3").Value = 'Pomelo'
$ws.Range("K3").Value = 'Start Ruby'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 430
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7500
$ws.Range("P3").Value = 7326
$ws.Range("Q3").Value = '$/caja 14 kilos granel'
$ws.Range("R3").Value = 'Región Metropolitana'
$ws.Range("S3").Value = 523
$ws.Range("T3").Value = 14
